$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row 2: drop trailing "*" from E2, add two new period columns
$ws.Range("E2").Value = "Apr 2023- Mar 2024"
$ws.Range("F2").Value = "Oct 2023- Sep 2024"
$ws.Range("G2").Value = "Apr 2024- Mar 2025*"

# Update data row 3: new/changed figures
$ws.Range("E3").Value = 560
$ws.Range("F3").Value = 560
$ws.Range("G3").Value = 539

# Copy formatting from the existing D column (header+data) to the new F/G columns
$ws.Range("E2").Copy()
$ws.Range("F2:G2").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("F3:G3").PasteSpecial(-4122)

# Update selection to mirror the saved workbook state
$ws.Range("D8").Select()
